$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '76.399.72'
$ws.Range("E2").Value = '  +0.00%  '
$ws.Range("D3").Value = '3.046.19'
$ws.Range("E3").Value = '  +3.58%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '198.62'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.08%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '618.36'
$ws.Range("D6").Style = "Normal"
$ws.Range("E8").Value = '  -0.85%  '
$ws.Range("E9").Value = '  +4.66%  '
$ws.Range("D10").Value = '3.046.02'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.439'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.93%  '
$ws.Range("E12").Value = '  -0.58%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.26'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +6.08%  '
$ws.Range("D14").Value = '3.607.65'
$ws.Range("E14").Value = '  +3.66%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '28.84'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +2.61%  '
$ws.Range("D16").Value = '76.309.60'
$ws.Range("E16").Value = '  -0.02%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000193'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.76%  '
$ws.Range("D18").Value = '3.048.60'
$ws.Range("E18").Value = '  +4.16%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.54'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.50%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '8.99'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +2.96%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '381.76'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +2.19%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '2.40'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +4.13%  '
$ws.Range("E23").Value = '  +1.00%  '
$ws.Range("D24").Value = '3.198.37'
$ws.Range("E24").Value = '  +4.60%  '
$ws.Range("E26").Value = '  +0.01%  '
$ws.Range("E27").Value = '  +1.14%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.78'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.81%  '
$ws.Range("E29").Value = '  -0.24%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.00'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.16%  '
$ws.Range("E31").Value = '  +4.64%  '
$ws.Range("E32").Value = '  +1.06%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '495.17'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.78%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.91'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +4.27%  '
$ws.Range("E35").Value = '  -0.08%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '20.59'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.58%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '163.23'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.30%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.117'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +5.07%  '
$ws.Range("E39").Value = '  +1.92%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '192.15'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +7.46%  '
$ws.Range("E41").Value = '  -4.46%  '
$ws.Range("E42").Value = '  -4.04%  '
$ws.Range("E43").Value = '  +0.04%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.799'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +20.72%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '5.11'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +3.17%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '41.91'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +4.26%  '
$ws.Range("E47").Value = '  +4.20%  '
$ws.Range("E48").Value = '  -0.81%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.43'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +4.46%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.600'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +2.48%  '
$ws.Range("E51").Value = '  -0.53%  '
